$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 63 (above the current "Chino / Primera / 44168" row),
# which pushes the existing rows 63-97 down to 64-98 and extends the used
# range from A1:R97 to A1:R98.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with this week's new price observation.
# All the "constant" columns (market/region/category/unit/origin/etc.) follow
# the same pattern as every other row in this sheet.
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = "Vega Monumental Concepción"
$ws.Range("C63").Value = "Bíobío"
$ws.Range("D63").Value = 44488
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = 100112003
$ws.Range("G63").Value = "Ajo"
$ws.Range("H63").Value = "Chino"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 400
$ws.Range("K63").Value = 16000
$ws.Range("L63").Value = 16500
$ws.Range("M63").Value = 16250
$ws.Range("N63").Value = "$/caja 10 kilos"
$ws.Range("O63").Value = "China"
$ws.Range("P63").Value = 1625
$ws.Range("Q63").Value = 10
$ws.Range("R63").Value = "Hortaliza"
